$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G4").Value = "fff"
$ws.Range("G15").Select()
